$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H92").Value = 456.56522
$ws_ALC.Range("I92").Value = 425.55
$ws_ALC.Range("J92").Value = 663.3333
$ws_ALC.Range("K92").Value = 425.55
$ws_ALC.Range("L92").Value = 663.3333
$ws_ALC.Range("M92").Value = 822.45
$ws_ALC.Range("N92").Value = -3159.3333
$ws_ALC.Range("H100").Value = 11496147
$ws_ALC.Range("I100").Value = 18519992
$ws_ALC.Range("J100").Value = 2582.6365
$ws_ALC.Range("K100").Value = 18519992
$ws_ALC.Range("L100").Value = 2582.6365
$ws_ALC.Range("M100").Value = -18519451
$ws_ALC.Range("N100").Value = -3664.6365
$ws_ALC.Range("H113").Value = 2172.45
$ws_ALC.Range("I113").Value = 1662.3077
$ws_ALC.Range("J113").Value = 3119.8572
$ws_ALC.Range("K113").Value = 1662.3077
$ws_ALC.Range("L113").Value = 3119.8572
$ws_ALC.Range("M113").Value = 1591.6923
$ws_ALC.Range("N113").Value = -9627.8572
$ws_ALC.Range("H132").Value = 20002146
$ws_ALC.Range("I132").Value = 26317308
$ws_ALC.Range("J132").Value = 4132.25
$ws_ALC.Range("K132").Value = 78951924
$ws_ALC.Range("L132").Value = 12396.75
$ws_ALC.Range("M132").Value = -78949394
$ws_ALC.Range("N132").Value = -17456.75
$ws_ALC.Range("H137").Value = 3963.0625
$ws_ALC.Range("I137").Value = 6946.3125
$ws_ALC.Range("J137").Value = 979.8125
$ws_ALC.Range("K137").Value = 20838.9375
$ws_ALC.Range("L137").Value = 2939.4375
$ws_ALC.Range("M137").Value = -18288.9375
$ws_ALC.Range("N137").Value = -8039.4375
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1218.9117
$ws_ARM.Range("I61").Value = 833.72
$ws_ARM.Range("J61").Value = 2288.889
$ws_ARM.Range("K61").Value = 833.72
$ws_ARM.Range("L61").Value = 2288.889
$ws_ARM.Range("M61").Value = -621.72
$ws_ARM.Range("N61").Value = -2712.889
$ws_ARM.Range("H74").Value = 501111.34
$ws_ARM.Range("I74").Value = 556656.4
$ws_ARM.Range("J74").Value = 1206
$ws_ARM.Range("K74").Value = 556656.4
$ws_ARM.Range("L74").Value = 1206
$ws_ARM.Range("M74").Value = -555782.4
$ws_ARM.Range("N74").Value = -2954
$ws_ARM.Range("H77").Value = 501111.34
$ws_ARM.Range("I77").Value = 556656.4
$ws_ARM.Range("J77").Value = 1206
$ws_ARM.Range("K77").Value = 2783282
$ws_ARM.Range("L77").Value = 6030
$ws_ARM.Range("M77").Value = -2778914
$ws_ARM.Range("N77").Value = -14766
$ws_ARM.Range("H136").Value = 1218.9117
$ws_ARM.Range("I136").Value = 833.72
$ws_ARM.Range("J136").Value = 2288.889
$ws_ARM.Range("K136").Value = 2501.16
$ws_ARM.Range("L136").Value = 6866.667
$ws_ARM.Range("M136").Value = 48.84000000000015
$ws_ARM.Range("N136").Value = -11966.667
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H94").Value = 1359.4117
$ws_BSM.Range("I94").Value = 1536.3636
$ws_BSM.Range("J94").Value = 1035
$ws_BSM.Range("K94").Value = 1536.3636
$ws_BSM.Range("L94").Value = 1035
$ws_BSM.Range("M94").Value = -1085.3636
$ws_BSM.Range("N94").Value = -1937
$ws_BSM.Range("H134").Value = 951.807
$ws_BSM.Range("I134").Value = 653.875
$ws_BSM.Range("J134").Value = 1652.8235
$ws_BSM.Range("K134").Value = 1961.625
$ws_BSM.Range("L134").Value = 4958.470499999999
$ws_BSM.Range("M134").Value = 573.375
$ws_BSM.Range("N134").Value = -10028.4705
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 7577668
$ws_CRP.Range("I31").Value = 10418175
$ws_CRP.Range("J31").Value = 2983.1667
$ws_CRP.Range("K31").Value = 10418175
$ws_CRP.Range("L31").Value = 2983.1667
$ws_CRP.Range("M31").Value = -10417880
$ws_CRP.Range("N31").Value = -3573.1667
$ws_CRP.Range("H34").Value = 7577668
$ws_CRP.Range("I34").Value = 10418175
$ws_CRP.Range("J34").Value = 2983.1667
$ws_CRP.Range("K34").Value = 10418175
$ws_CRP.Range("L34").Value = 2983.1667
$ws_CRP.Range("M34").Value = -10417973
$ws_CRP.Range("N34").Value = -3387.1667
$ws_CRP.Range("H58").Value = 1263.3448
$ws_CRP.Range("I58").Value = 971.9524
$ws_CRP.Range("J58").Value = 2028.25
$ws_CRP.Range("K58").Value = 971.9524
$ws_CRP.Range("L58").Value = 2028.25
$ws_CRP.Range("M58").Value = -768.9524
$ws_CRP.Range("N58").Value = -2434.25
$ws_CRP.Range("H132").Value = 26531.64
$ws_CRP.Range("I132").Value = 823.36365
$ws_CRP.Range("J132").Value = 167927.17
$ws_CRP.Range("K132").Value = 2470.09095
$ws_CRP.Range("L132").Value = 503781.51
$ws_CRP.Range("M132").Value = 59.90905000000021
$ws_CRP.Range("N132").Value = -508841.51
$ws_CRP.Range("H134").Value = 17514.35
$ws_CRP.Range("I134").Value = 21736.041
$ws_CRP.Range("J134").Value = 627.5833
$ws_CRP.Range("K134").Value = 65208.12300000001
$ws_CRP.Range("L134").Value = 1882.7499
$ws_CRP.Range("M134").Value = -62673.12300000001
$ws_CRP.Range("N134").Value = -6952.7499
$ws_CRP.Range("H136").Value = 1263.3448
$ws_CRP.Range("I136").Value = 971.9524
$ws_CRP.Range("J136").Value = 2028.25
$ws_CRP.Range("K136").Value = 2915.8572
$ws_CRP.Range("L136").Value = 6084.75
$ws_CRP.Range("M136").Value = -365.8571999999999
$ws_CRP.Range("N136").Value = -11184.75
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 695.93335
$ws_CUL.Range("I5").Value = 695.93335
$ws_CUL.Range("J5").Value = 0
$ws_CUL.Range("K5").Value = 2087.80005
$ws_CUL.Range("L5").Value = 0
$ws_CUL.Range("M5").Value = -1975.80005
$ws_CUL.Range("N5").ClearContents()
$ws_CUL.Range("H122").Value = 871.1429000000001
$ws_CUL.Range("I122").Value = 699.6667
$ws_CUL.Range("J122").Value = 999.75
$ws_CUL.Range("K122").Value = 6297.0003
$ws_CUL.Range("L122").Value = 8997.75
$ws_CUL.Range("M122").Value = -3847.0003
$ws_CUL.Range("N122").Value = -13897.75
$ws_CUL.Range("H131").Value = 6579941.5
$ws_CUL.Range("I131").Value = 871.5833
$ws_CUL.Range("J131").Value = 7813517
$ws_CUL.Range("K131").Value = 2614.7499
$ws_CUL.Range("L131").Value = 23440551
$ws_CUL.Range("M131").Value = 2425.2501
$ws_CUL.Range("N131").Value = -23450631
$ws_CUL.Range("H135").Value = 695.93335
$ws_CUL.Range("I135").Value = 695.93335
$ws_CUL.Range("J135").Value = 0
$ws_CUL.Range("K135").Value = 6263.40015
$ws_CUL.Range("L135").Value = 0
$ws_CUL.Range("M135").Value = -3728.40015
$ws_CUL.Range("N135").ClearContents()
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 22750.447
$ws_GSM.Range("I132").Value = 30846.766
$ws_GSM.Range("J132").Value = 1575.4615
$ws_GSM.Range("K132").Value = 92540.298
$ws_GSM.Range("L132").Value = 4726.3845
$ws_GSM.Range("M132").Value = -90010.298
$ws_GSM.Range("N132").Value = -9786.3845
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H40").Value = 1994.7826
$ws_LTW.Range("I40").Value = 1825
$ws_LTW.Range("J40").Value = 2180
$ws_LTW.Range("K40").Value = 1825
$ws_LTW.Range("L40").Value = 2180
$ws_LTW.Range("M40").Value = -1689
$ws_LTW.Range("N40").Value = -2452
$ws_LTW.Range("H50").Value = 7000
$ws_LTW.Range("J50").Value = 7000
$ws_LTW.Range("L50").Value = 7000
$ws_LTW.Range("N50").Value = -8274
$ws_LTW.Range("H82").Value = 1946.2727
$ws_LTW.Range("I82").Value = 2022.7142
$ws_LTW.Range("J82").Value = 1812.5
$ws_LTW.Range("K82").Value = 2022.7142
$ws_LTW.Range("L82").Value = 1812.5
$ws_LTW.Range("M82").Value = -1661.7142
$ws_LTW.Range("N82").Value = -2534.5
$ws_LTW.Range("H85").Value = 1946.2727
$ws_LTW.Range("I85").Value = 2022.7142
$ws_LTW.Range("J85").Value = 1812.5
$ws_LTW.Range("K85").Value = 2022.7142
$ws_LTW.Range("L85").Value = 1812.5
$ws_LTW.Range("M85").Value = -774.7141999999999
$ws_LTW.Range("N85").Value = -4308.5
$ws_LTW.Range("H132").Value = 2997.0466
$ws_LTW.Range("I132").Value = 3228.7646
$ws_LTW.Range("J132").Value = 2121.6667
$ws_LTW.Range("K132").Value = 9686.293799999999
$ws_LTW.Range("L132").Value = 6365.000100000001
$ws_LTW.Range("M132").Value = -7156.293799999999
$ws_LTW.Range("N132").Value = -11425.0001
$ws_LTW.Range("H136").Value = 2145.9607
$ws_LTW.Range("I136").Value = 1331.6364
$ws_LTW.Range("J136").Value = 3638.889
$ws_LTW.Range("K136").Value = 3994.9092
$ws_LTW.Range("L136").Value = 10916.667
$ws_LTW.Range("M136").Value = -1444.9092
$ws_LTW.Range("N136").Value = -16016.667
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H100").Value = 1158.25
$ws_WVR.Range("I100").Value = 1427.2222
$ws_WVR.Range("J100").Value = 996.86664
$ws_WVR.Range("K100").Value = 2854.4444
$ws_WVR.Range("L100").Value = 1993.73328
$ws_WVR.Range("M100").Value = -2313.4444
$ws_WVR.Range("N100").Value = -3075.73328
$ws_WVR.Range("H132").Value = 18942044
$ws_WVR.Range("I132").Value = 27175128
$ws_WVR.Range("J132").Value = 5952.25
$ws_WVR.Range("K132").Value = 81525384
$ws_WVR.Range("L132").Value = 17856.75
$ws_WVR.Range("M132").Value = -81522854
$ws_WVR.Range("N132").Value = -22916.75
